# Weekly price-list update: insert one new daily record ahead of the
# existing "Arveja Verde" (Vega Central Mapocho de Santiago) rows.
#
# The new record becomes row 54; every existing row from 54 downward
# (through the former last row, 83) shifts down by one to make room,
# ending at row 84 — matching the new dimension A1:R84.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data (old rows 54:83) down by one row.
$ws.Rows("54:54").Insert()

# Populate the newly inserted row 54 with the new weekly observation.
$ws.Range("A54").Value = 9
$ws.Range("B54").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C54").Value = "Metropolitana"
$ws.Range("D54").Value = 44529
$ws.Range("E54").Value = 13
$ws.Range("F54").Value = 100112022
$ws.Range("G54").Value = "Arveja Verde"
$ws.Range("H54").Value = "Sin especificar"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 34
$ws.Range("K54").Value = 16000
$ws.Range("L54").Value = 18000
$ws.Range("M54").Value = 17000
$ws.Range("N54").Value = "$/saco 25 kilos"
$ws.Range("O54").Value = "Carahue"
$ws.Range("P54").Value = 680
$ws.Range("Q54").Value = 25
$ws.Range("R54").Value = "Hortaliza"
